# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.079.13"
$ws.Range("E2").Value = "  +4.36%  "
$ws.Range("D3").Value = "'1.906.70"
$ws.Range("E3").Value = "  +5.40%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'251.75"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.5103"
$ws.Range("E7").Value = "  +2.77%  "
$ws.Range("D8").Value = "'45.13"
$ws.Range("E8").Value = "  +4.40%  "
$ws.Range("D9").Value = "'0.3016"
$ws.Range("E9").Value = "  +8.61%  "
$ws.Range("D10").Value = "'0.06800"
$ws.Range("E10").Value = "  +6.05%  "
$ws.Range("D11").Value = "'1.907.55"
$ws.Range("E11").Value = "  +5.44%  "
$ws.Range("D12").Value = "'17.27"
$ws.Range("E12").Value = "  +3.20%  "
$ws.Range("D13").Value = "'0.07330"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("D14").Value = "'0.6979"
$ws.Range("E14").Value = "  +8.01%  "
$ws.Range("D15").Value = "'86.59"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("D16").Value = "'4.910"
$ws.Range("E16").Value = "  +4.74%  "
$ws.Range("D17").Value = "'30.068.25"
$ws.Range("E17").Value = "  +4.42%  "
$ws.Range("D18").Value = "'0.000008188"
$ws.Range("E18").Value = "  +11.41%  "
$ws.Range("D19").Value = "'0.9999"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'13.04"
$ws.Range("E20").Value = "  +6.57%  "
$ws.Range("D21").Value = "'2.154.29"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'4.830"
$ws.Range("E23").Value = "  +5.64%  "
$ws.Range("D24").Value = "'5.736"
$ws.Range("E24").Value = "  +7.42%  "
$ws.Range("D25").Value = "'9.271"
$ws.Range("E25").Value = "  +4.23%  "
$ws.Range("D26").Value = "'147.68"
$ws.Range("E26").Value = "  +3.75%  "
$ws.Range("D27").Value = "'135.00"
$ws.Range("E27").Value = "  +4.37%  "
$ws.Range("D28").Value = "'17.06"
$ws.Range("E28").Value = "  +4.03%  "
$ws.Range("D29").Value = "'1.995"
$ws.Range("E29").Value = "  +6.07%  "
$ws.Range("D30").Value = "'1.408"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").Value = "'4.256"
$ws.Range("E31").Value = "  +3.08%  "
$ws.Range("D32").Value = "'0.08813"
$ws.Range("E32").Value = "  +5.70%  "
$ws.Range("D33").Value = "'3.999"
$ws.Range("E33").Value = "  +5.20%  "
$ws.Range("D34").Value = "'0.05052"
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D35").Value = "'1.141"
$ws.Range("D36").Value = "'0.7210"
$ws.Range("E36").Value = "  +7.86%  "
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").Value = "'2.814"
$ws.Range("E38").Value = "  +3.24%  "
$ws.Range("D39").Value = "'2.267"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'0.9641"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").Value = "'0.01695"
$ws.Range("E41").Value = "  +6.74%  "
$ws.Range("D42").Value = "'6.144"
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").Value = "'0.4312"
$ws.Range("E43").Value = "  +5.91%  "
$ws.Range("D44").Value = "'104.72"
$ws.Range("E44").Value = "  +5.45%  "
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "'7.617"
$ws.Range("E46").Value = "  +6.46%  "
$ws.Range("D47").Value = "'0.1280"
$ws.Range("E47").Value = "  +5.04%  "
$ws.Range("D48").Value = "'0.05741"
$ws.Range("E48").Value = "  +4.12%  "
$ws.Range("E49").Value = "  +5.91%  "
$ws.Range("D50").Value = "'8.426"
$ws.Range("E50").Value = "  +3.67%  "
$ws.Range("D51").Value = "'0.3816"
$ws.Range("E51").Value = "  +5.36%  "
